$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '28.919.05'
Set-TextValue $ws 'E2' '  -1.69%  '

Set-TextValue $ws 'D3' '1.829.41'
Set-TextValue $ws 'E3' '  -1.94%  '

Set-TextValue $ws 'E4' '  -0.13%  '

Set-TextValue $ws 'D5' '240.50'
Set-TextValue $ws 'E5' '  -1.11%  '

Set-TextValue $ws 'D6' '0.6875'
Set-TextValue $ws 'E6' '  -2.80%  '

Set-TextValue $ws 'D7' '0.9997'
Set-TextValue $ws 'E7' '  -0.10%  '

Set-TextValue $ws 'B8' 'Dogecoin'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws 'D8' '0.07641'
Set-TextValue $ws 'E8' '  -3.01%  '

Set-TextValue $ws 'B9' 'Cardano'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws 'D9' '0.3044'
Set-TextValue $ws 'E9' '  -2.70%  '

Set-TextValue $ws 'D10' '23.54'
Set-TextValue $ws 'E10' '  -3.78%  '

Set-TextValue $ws 'D11' '0.07780'
Set-TextValue $ws 'E11' '  -2.72%  '

Set-TextValue $ws 'D12' '1.825.57'
Set-TextValue $ws 'E12' '  -2.24%  '

Set-TextValue $ws 'D13' '5.062'
Set-TextValue $ws 'E13' '  -2.56%  '

Set-TextValue $ws 'D14' '90.44'
Set-TextValue $ws 'E14' '  -3.28%  '

Set-TextValue $ws 'D15' '0.6746'
Set-TextValue $ws 'E15' '  -3.50%  '

Set-TextValue $ws 'D16' '6.428'
Set-TextValue $ws 'E16' '  -0.21%  '

Set-TextValue $ws 'D17' '0.000008271'
Set-TextValue $ws 'E17' '  -1.11%  '

Set-TextValue $ws 'D18' '28.900.15'
Set-TextValue $ws 'E18' '  -1.74%  '

Set-TextValue $ws 'D19' '242.50'
Set-TextValue $ws 'E19' '  -4.16%  '

Set-TextValue $ws 'D20' '2.074.59'
Set-TextValue $ws 'E20' '  -2.14%  '

Set-TextValue $ws 'D21' '12.66'
Set-TextValue $ws 'E21' '  -3.35%  '

Set-TextValue $ws 'D22' '0.9999'
Set-TextValue $ws 'E22' '  -0.10%  '

Set-TextValue $ws 'D23' '7.413'
Set-TextValue $ws 'E23' '  -2.70%  '

Set-TextValue $ws 'E24' '  -0.16%  '

Set-TextValue $ws 'D25' '0.1474'
Set-TextValue $ws 'E25' '  -5.59%  '

Set-TextValue $ws 'D26' '161.36'
Set-TextValue $ws 'E26' '  +0.39%  '

Set-TextValue $ws 'D27' '8.774'
Set-TextValue $ws 'E27' '  -2.55%  '

Set-TextValue $ws 'D28' '18.19'

Set-TextValue $ws 'D29' '1.535'
Set-TextValue $ws 'E29' '  +2.45%  '

Set-TextValue $ws 'D30' '4.209'
Set-TextValue $ws 'E30' '  -2.62%  '

Set-TextValue $ws 'D31' '4.129'
Set-TextValue $ws 'E31' '  -3.54%  '

Set-TextValue $ws 'D32' '1.190'
Set-TextValue $ws 'E32' '  -1.72%  '

Set-TextValue $ws 'D33' '0.05110'
Set-TextValue $ws 'E33' '  -3.79%  '

Set-TextValue $ws 'D34' '0.7472'
Set-TextValue $ws 'E34' '  -0.55%  '

Set-TextValue $ws 'D35' '1.820'
Set-TextValue $ws 'E35' '  -3.53%  '

Set-TextValue $ws 'D36' '1.142'
Set-TextValue $ws 'E36' '  -2.18%  '

Set-TextValue $ws 'D37' '2.677'
Set-TextValue $ws 'E37' '  -1.24%  '

Set-TextValue $ws 'D38' '0.01840'
Set-TextValue $ws 'E38' '  -2.05%  '

Set-TextValue $ws 'D39' '1.216.24'
Set-TextValue $ws 'E39' '  -4.76%  '

Set-TextValue $ws 'D40' '2.679'
Set-TextValue $ws 'E40' '  -2.23%  '

Set-TextValue $ws 'D41' '0.9160'
Set-TextValue $ws 'E41' '  +2.02%  '

Set-TextValue $ws 'D42' '108.36'
Set-TextValue $ws 'E42' '  -0.43%  '

Set-TextValue $ws 'D43' '0.9994'
Set-TextValue $ws 'E43' '  -0.13%  '

Set-TextValue $ws 'D44' '0.5168'
Set-TextValue $ws 'E44' '  -0.09%  '

Set-TextValue $ws 'D45' '1.974.26'
Set-TextValue $ws 'E45' '  -2.47%  '

Set-TextValue $ws 'D46' '9.496'
Set-TextValue $ws 'E46' '  -1.11%  '

Set-TextValue $ws 'E47' '  -5.28%  '

Set-TextValue $ws 'D48' '5.328'
Set-TextValue $ws 'E48' '  -11.59%  '

Set-TextValue $ws 'B49' 'Aave'
Set-TextValue $ws 'C49' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws 'D49' '62.97'
Set-TextValue $ws 'E49' '  -11.61%  '

Set-TextValue $ws 'B50' 'RenderToken'
Set-TextValue $ws 'C50' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws 'D50' '1.729'
Set-TextValue $ws 'E50' '  -3.45%  '

Set-TextValue $ws 'D51' '0.4183'
Set-TextValue $ws 'E51' '  -2.83%  '
